# "Generate Report for Handback"
#
# For each language sheet (zh-cn, de-de) the handback report now records,
# per source file, BOTH the already-existing "Latest Handoff" columns
# (C/D) AND a newly populated "Latest Target"/"Latest Handback File"
# pair (E/F), plus an updated "Latest Handback DateTime" (G) showing the
# file came back in sync with en-US. The Status column (B) is updated to
# reflect the handback, too.

$wb = $excel.ActiveWorkbook

$langSheets = @(
    @{ Name = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/016dadc0cc36b75d449fbd4718c6e2591b58e36c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf";
       HandbackDateTime = "2016-02-24 09:45:45" },
    @{ Name = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/381b3aeb097f5ac9ba76dfa127f865212072eaf6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf";
       HandbackDateTime = "2016-02-24 09:46:15" }
)

$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/cc8fa6c0f710810780268ba1c56e904adcfc861e/e2e/a.md"
$statusText = "Handed back: in sync with en-US"

# The "Overview" sheet mirrors the same per-language status in its
# zh-cn / de-de columns (B/C) for each source file row, so it picks up
# the same wording as the per-language sheets' Status column.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (B) for the two real source files (a.md / b.md rows).
    $ws.Range("B2").Value = $statusText
    $ws.Range("B3").Value = $statusText

    # Newly tracked "Latest Target File" (E) / "Latest Handback File" (F)
    # columns for rows 2 and 3 -- same file identities as the existing
    # handoff columns (A = source file, C = handoff xlf).
    $ws.Range("E2").Value = "a.md"
    $ws.Range("F2").Value = $lang.Xlf
    $ws.Range("E3").Value = "a.md"
    $ws.Range("F3").Value = $lang.Xlf

    # Latest Handback DateTime (G) moves from the "never handed back"
    # sentinel to the real handback timestamp for rows 2 and 3.
    $ws.Range("G2").Value = $lang.HandbackDateTime
    $ws.Range("G3").Value = $lang.HandbackDateTime

    # Mirror the existing hyperlink styling used on A2/A3 (-> a.md / b.md)
    # and C2/C3 (-> handoff xlf) onto the new E/F cells.
    $ws.Hyperlinks.Add($ws.Range("E2"), $aMdUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.HandoffUrl, "", "", $lang.Xlf)
    $ws.Hyperlinks.Add($ws.Range("E3"), $aMdUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("F3"), $lang.HandoffUrl, "", "", $lang.Xlf)
}
